$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE2").Value = 31.1
$ws.Range("AE3").Value = 84.40000000000001
$ws.Range("AE4").Value = 42.2
$ws.Range("AG5").Value = 93.3
$ws.Range("AE6").Value = 26.7
$ws.Range("AE7").Value = 31.1
$ws.Range("AE8").Value = 25.6
$ws.Range("AE9").Value = 40
$ws.Range("AG11").Value = 5
$ws.Range("AE12").Value = 62.2
$ws.Range("AG12").Value = 84.40000000000001
$ws.Range("AE13").Value = 40
$ws.Range("AG13").Value = 74.40000000000001
$ws.Range("AE14").Value = 5.6
$ws.Range("AE15").Value = 91.09999999999999
$ws.Range("AE16").Value = 33.3
$ws.Range("AG16").Value = 47.8
$ws.Range("AE17").Value = 86.7
$ws.Range("AG17").Value = 77.8
$ws.Range("AE18").Value = 21.1
$ws.Range("AE20").Value = 44.4
$ws.Range("AG20").Value = 50
$ws.Range("AE21").Value = 10
$ws.Range("AE22").Value = 53.3
$ws.Range("AG22").Value = 41.1
$ws.Range("AE23").Value = 60
$ws.Range("AE24").Value = 71.09999999999999
$ws.Range("AG24").Value = 37.8
$ws.Range("AE25").Value = 81.09999999999999
$ws.Range("AG25").Value = 88.90000000000001
$ws.Range("AE26").Value = 64.40000000000001
$ws.Range("AG26").Value = 71.09999999999999
$ws.Range("AE27").Value = 22.2
$ws.Range("AG27").Value = 34.4
$ws.Range("AG28").Value = 75.59999999999999
$ws.Range("AE30").Value = 50
$ws.Range("AG31").Value = 23.3
$ws.Range("AE32").Value = 84.40000000000001
$ws.Range("AG32").Value = 88.90000000000001
$ws.Range("AE33").Value = 13.3
$ws.Range("AG33").Value = 7.8
$ws.Range("AE34").Value = 73.3
$ws.Range("AE35").Value = 20
$ws.Range("AG35").Value = 5.6
$ws.Range("AE36").Value = 15.6
$ws.Range("AE38").Value = 66.7
$ws.Range("AE39").Value = 68.90000000000001
$ws.Range("AG40").Value = 94.40000000000001
$ws.Range("AE42").Value = 55.6
$ws.Range("AE43").Value = 37.8
$ws.Range("AG43").Value = 18.9
$ws.Range("AG44").Value = 45.6
$ws.Range("AE45").Value = 61.1
$ws.Range("AG45").Value = 38.9
$ws.Range("AE46").Value = 91.09999999999999
$ws.Range("AG46").Value = 85.59999999999999
